# Appends the new rows captured by the latest SeniorConnect sensor sweep
# (2026-02-01, ~18:27-18:28) to the PIR, Humidity and Temperature logs.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# PIR sheet - append new sensor-log rows
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("PIR")

# Row 46
$ws.Cells.Item(46,1).NumberFormat = "@"
$ws.Cells.Item(46,1).Value = "2026-02-01"
$ws.Cells.Item(46,2).Value = "18:27:39"
$ws.Cells.Item(46,3).Value = "18:00"
$ws.Cells.Item(46,4).Value = "Bathroom"
$ws.Cells.Item(46,5).Value = "No Motion"
$ws.Cells.Item(46,6).Value = "Inactive"

# Row 47
$ws.Cells.Item(47,1).NumberFormat = "@"
$ws.Cells.Item(47,1).Value = "2026-02-01"
$ws.Cells.Item(47,2).Value = "18:27:41"
$ws.Cells.Item(47,3).Value = "18:00"
$ws.Cells.Item(47,4).Value = "Bathroom"
$ws.Cells.Item(47,5).Value = "Motion Detected"
$ws.Cells.Item(47,6).Value = "Active"

# ---------------------------------------------------------------
# Humidity sheet - append new sensor-log rows
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Humidity")

# Row 94
$ws.Cells.Item(94,1).NumberFormat = "@"
$ws.Cells.Item(94,1).Value = "2026-02-01"
$ws.Cells.Item(94,2).Value = "18:27:35"
$ws.Cells.Item(94,3).Value = "18:00"
$ws.Cells.Item(94,4).Value = "Bathroom"
$ws.Cells.Item(94,5).NumberFormat = "@"
$ws.Cells.Item(94,5).Value = "81.9%"
$ws.Cells.Item(94,6).Value = "Active"

# Row 95
$ws.Cells.Item(95,1).NumberFormat = "@"
$ws.Cells.Item(95,1).Value = "2026-02-01"
$ws.Cells.Item(95,2).Value = "18:27:39"
$ws.Cells.Item(95,3).Value = "18:00"
$ws.Cells.Item(95,4).Value = "Bathroom"
$ws.Cells.Item(95,5).NumberFormat = "@"
$ws.Cells.Item(95,5).Value = "82.5%"
$ws.Cells.Item(95,6).Value = "Active"

# Row 96
$ws.Cells.Item(96,1).NumberFormat = "@"
$ws.Cells.Item(96,1).Value = "2026-02-01"
$ws.Cells.Item(96,2).Value = "18:27:43"
$ws.Cells.Item(96,3).Value = "18:00"
$ws.Cells.Item(96,4).Value = "Bathroom"
$ws.Cells.Item(96,5).NumberFormat = "@"
$ws.Cells.Item(96,5).Value = "79.4%"
$ws.Cells.Item(96,6).Value = "Active"

# Row 97
$ws.Cells.Item(97,1).NumberFormat = "@"
$ws.Cells.Item(97,1).Value = "2026-02-01"
$ws.Cells.Item(97,2).Value = "18:27:49"
$ws.Cells.Item(97,3).Value = "18:00"
$ws.Cells.Item(97,4).Value = "Bathroom"
$ws.Cells.Item(97,5).NumberFormat = "@"
$ws.Cells.Item(97,5).Value = "80.0%"
$ws.Cells.Item(97,6).Value = "Active"

# Row 98
$ws.Cells.Item(98,1).NumberFormat = "@"
$ws.Cells.Item(98,1).Value = "2026-02-01"
$ws.Cells.Item(98,2).Value = "18:27:54"
$ws.Cells.Item(98,3).Value = "18:00"
$ws.Cells.Item(98,4).Value = "Bathroom"
$ws.Cells.Item(98,5).NumberFormat = "@"
$ws.Cells.Item(98,5).Value = "78.3%"
$ws.Cells.Item(98,6).Value = "Active"

# Row 99
$ws.Cells.Item(99,1).NumberFormat = "@"
$ws.Cells.Item(99,1).Value = "2026-02-01"
$ws.Cells.Item(99,2).Value = "18:28:09"
$ws.Cells.Item(99,3).Value = "18:00"
$ws.Cells.Item(99,4).Value = "Bathroom"
$ws.Cells.Item(99,5).NumberFormat = "@"
$ws.Cells.Item(99,5).Value = "78.9%"
$ws.Cells.Item(99,6).Value = "Active"

# Row 100
$ws.Cells.Item(100,1).NumberFormat = "@"
$ws.Cells.Item(100,1).Value = "2026-02-01"
$ws.Cells.Item(100,2).Value = "18:28:14"
$ws.Cells.Item(100,3).Value = "18:00"
$ws.Cells.Item(100,4).Value = "Bathroom"
$ws.Cells.Item(100,5).NumberFormat = "@"
$ws.Cells.Item(100,5).Value = "78.0%"
$ws.Cells.Item(100,6).Value = "Active"

# Row 101
$ws.Cells.Item(101,1).NumberFormat = "@"
$ws.Cells.Item(101,1).Value = "2026-02-01"
$ws.Cells.Item(101,2).Value = "18:28:19"
$ws.Cells.Item(101,3).Value = "18:00"
$ws.Cells.Item(101,4).Value = "Bathroom"
$ws.Cells.Item(101,5).NumberFormat = "@"
$ws.Cells.Item(101,5).Value = "79.0%"
$ws.Cells.Item(101,6).Value = "Active"

# Row 102
$ws.Cells.Item(102,1).NumberFormat = "@"
$ws.Cells.Item(102,1).Value = "2026-02-01"
$ws.Cells.Item(102,2).Value = "18:28:24"
$ws.Cells.Item(102,3).Value = "18:00"
$ws.Cells.Item(102,4).Value = "Bathroom"
$ws.Cells.Item(102,5).NumberFormat = "@"
$ws.Cells.Item(102,5).Value = "78.2%"
$ws.Cells.Item(102,6).Value = "Active"

# Row 103
$ws.Cells.Item(103,1).NumberFormat = "@"
$ws.Cells.Item(103,1).Value = "2026-02-01"
$ws.Cells.Item(103,2).Value = "18:28:29"
$ws.Cells.Item(103,3).Value = "18:00"
$ws.Cells.Item(103,4).Value = "Bathroom"
$ws.Cells.Item(103,5).NumberFormat = "@"
$ws.Cells.Item(103,5).Value = "79.1%"
$ws.Cells.Item(103,6).Value = "Active"

# ---------------------------------------------------------------
# Temperature sheet - append new sensor-log rows
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Temperature")

# Row 94
$ws.Cells.Item(94,1).NumberFormat = "@"
$ws.Cells.Item(94,1).Value = "2026-02-01"
$ws.Cells.Item(94,2).Value = "18:27:35"
$ws.Cells.Item(94,3).Value = "18:00"
$ws.Cells.Item(94,4).Value = "Bathroom"
$ws.Cells.Item(94,5).Value = "29.9C"
$ws.Cells.Item(94,6).Value = "Active"

# Row 95
$ws.Cells.Item(95,1).NumberFormat = "@"
$ws.Cells.Item(95,1).Value = "2026-02-01"
$ws.Cells.Item(95,2).Value = "18:27:39"
$ws.Cells.Item(95,3).Value = "18:00"
$ws.Cells.Item(95,4).Value = "Bathroom"
$ws.Cells.Item(95,5).Value = "29.9C"
$ws.Cells.Item(95,6).Value = "Active"

# Row 96
$ws.Cells.Item(96,1).NumberFormat = "@"
$ws.Cells.Item(96,1).Value = "2026-02-01"
$ws.Cells.Item(96,2).Value = "18:27:44"
$ws.Cells.Item(96,3).Value = "18:00"
$ws.Cells.Item(96,4).Value = "Bathroom"
$ws.Cells.Item(96,5).Value = "29.9C"
$ws.Cells.Item(96,6).Value = "Active"

# Row 97
$ws.Cells.Item(97,1).NumberFormat = "@"
$ws.Cells.Item(97,1).Value = "2026-02-01"
$ws.Cells.Item(97,2).Value = "18:27:49"
$ws.Cells.Item(97,3).Value = "18:00"
$ws.Cells.Item(97,4).Value = "Bathroom"
$ws.Cells.Item(97,5).Value = "29.8C"
$ws.Cells.Item(97,6).Value = "Active"

# Row 98
$ws.Cells.Item(98,1).NumberFormat = "@"
$ws.Cells.Item(98,1).Value = "2026-02-01"
$ws.Cells.Item(98,2).Value = "18:27:54"
$ws.Cells.Item(98,3).Value = "18:00"
$ws.Cells.Item(98,4).Value = "Bathroom"
$ws.Cells.Item(98,5).Value = "29.9C"
$ws.Cells.Item(98,6).Value = "Active"

# Row 99
$ws.Cells.Item(99,1).NumberFormat = "@"
$ws.Cells.Item(99,1).Value = "2026-02-01"
$ws.Cells.Item(99,2).Value = "18:28:09"
$ws.Cells.Item(99,3).Value = "18:00"
$ws.Cells.Item(99,4).Value = "Bathroom"
$ws.Cells.Item(99,5).Value = "29.8C"
$ws.Cells.Item(99,6).Value = "Active"

# Row 100
$ws.Cells.Item(100,1).NumberFormat = "@"
$ws.Cells.Item(100,1).Value = "2026-02-01"
$ws.Cells.Item(100,2).Value = "18:28:14"
$ws.Cells.Item(100,3).Value = "18:00"
$ws.Cells.Item(100,4).Value = "Bathroom"
$ws.Cells.Item(100,5).Value = "29.7C"
$ws.Cells.Item(100,6).Value = "Active"

# Row 101
$ws.Cells.Item(101,1).NumberFormat = "@"
$ws.Cells.Item(101,1).Value = "2026-02-01"
$ws.Cells.Item(101,2).Value = "18:28:19"
$ws.Cells.Item(101,3).Value = "18:00"
$ws.Cells.Item(101,4).Value = "Bathroom"
$ws.Cells.Item(101,5).Value = "29.7C"
$ws.Cells.Item(101,6).Value = "Active"

# Row 102
$ws.Cells.Item(102,1).NumberFormat = "@"
$ws.Cells.Item(102,1).Value = "2026-02-01"
$ws.Cells.Item(102,2).Value = "18:28:24"
$ws.Cells.Item(102,3).Value = "18:00"
$ws.Cells.Item(102,4).Value = "Bathroom"
$ws.Cells.Item(102,5).Value = "29.7C"
$ws.Cells.Item(102,6).Value = "Active"

# Row 103
$ws.Cells.Item(103,1).NumberFormat = "@"
$ws.Cells.Item(103,1).Value = "2026-02-01"
$ws.Cells.Item(103,2).Value = "18:28:29"
$ws.Cells.Item(103,3).Value = "18:00"
$ws.Cells.Item(103,4).Value = "Bathroom"
$ws.Cells.Item(103,5).Value = "29.6C"
$ws.Cells.Item(103,6).Value = "Active"
